# Update countries & provincias Spain
# - Refresh COVID-19 counters for several countries.
# - Because some case counts overtook neighboring countries, the ranking
#   (the sheet is sorted descending by column B "Casos totales") changes
#   for a few adjacent rows, which swaps which country name sits on which
#   row.
# - Update the "Datos actualizados..." timestamp note in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

function Set-Row($Row, $Name, $B, $C, $D, $E, $F, $G, $H) {
    $ws.Cells.Item($Row, 1).Value = $Name
    $ws.Cells.Item($Row, 2).Value = $B
    $ws.Cells.Item($Row, 3).Value = $C
    $ws.Cells.Item($Row, 4).Value = $D
    $ws.Cells.Item($Row, 5).Value = $E
    $ws.Cells.Item($Row, 6).Value = $F
    $ws.Cells.Item($Row, 7).Value = $G
    $ws.Cells.Item($Row, 8).Value = $H
}

# Row 15: Pakistan - values refreshed, ranking unchanged.
Set-Row 15 "Pakistan" 202955 4072 92624 106213 0 83 4118

# Rows 74/75: Uzbekistan overtakes Australia (7725 > 7686).
Set-Row 74 "Uzbekistan" 7725 43 5240 2465 0 0 20
Set-Row 75 "Australia" 7686 45 6979 603 0 0 104

# Rows 88/89/90: Kirguistan overtakes Republica de Yibuti and Bulgaria.
Set-Row 88 "Kirguistan" 4748 235 2242 2459 0 1 47
Set-Row 89 "Republica de Yibuti" 4643 0 4348 243 0 0 52
Set-Row 90 "Bulgaria" 4625 0 2475 1934 0 0 216

# Row 186: Butan - values refreshed, ranking unchanged.
Set-Row 186 "Butan" 76 1 38 38 0 0 0

# Rows 202/203: Laos and Santa Lucia tie (both 19) - swap order.
Set-Row 202 "Laos" 19 0 19 0 0 0 0
Set-Row 203 "Santa Lucia" 19 0 19 0 0 0 0

# Rows 209/210: Groenlandia and Islas Malvinas tie (both 13) - swap order.
Set-Row 209 "Groenlandia" 13 0 13 0 0 0 0
Set-Row 210 "Islas Malvinas" 13 0 13 0 0 0 0

# Update the "last updated" note timestamp.
$ws.Range("A1").Value = "Datos actualizados a 28 de Junio de 2020 a las 07:38"
